# "Changes of 5th May 2022"
# The three test-result cells in column R (rows 2-4) flip from "FAIL" to
# "PASS" (a new shared string "PASS" is introduced by this).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = "PASS"
$ws.Range("R3").Value = "PASS"
$ws.Range("R4").Value = "PASS"
